# Auto-generated Excel COM-interop script to apply cryptos.xlsx update
# Updates coin list ranking shifts (rows 6-18) and refreshed price/volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'311.91"
$ws.Range("E2").Value = "'0.93%"

# Row 3
$ws.Range("D3").Value = "'37.67"
$ws.Range("E3").Value = "'-0.06%"

# Row 4
$ws.Range("D4").Value = "'5.123"
$ws.Range("E4").Value = "'0.34%"

# Row 5
$ws.Range("D5").Value = "'0.07897"
$ws.Range("E5").Value = "'0.62%"

# Row 6
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.410"
$ws.Range("E6").Value = "'1.08%"

# Row 7
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.903"
$ws.Range("E7").Value = "'-3.53%"

# Row 8
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'8.258"
$ws.Range("E8").Value = "'-0.44%"

# Row 9
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.894"
$ws.Range("E9").Value = "'-8.25%"

# Row 10
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "'0.9286"
$ws.Range("E10").Value = "'-0.08%"

# Row 11
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.1223"
$ws.Range("E11").Value = "'-9.37%"

# Row 12
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1923"
$ws.Range("E12").Value = "'-3.83%"

# Row 13
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.09095"
$ws.Range("E13").Value = "'0.87%"

# Row 14
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03330"
$ws.Range("E14").Value = "'-3.07%"

# Row 15
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09632"
$ws.Range("E15").Value = "'-0.98%"

# Row 16
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001379"
$ws.Range("E16").Value = "'-1.14%"

# Row 17
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005750"
$ws.Range("E17").Value = "'-2.53%"

# Row 18
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.539"
$ws.Range("E18").Value = "'-1.26%"

# Row 19
$ws.Range("E19").Value = "'-1.68%"

# Row 20
$ws.Range("D20").Value = "'5.304"
$ws.Range("E20").Value = "'5.85%"

# Row 21
$ws.Range("D21").Value = "'0.1279"
$ws.Range("E21").Value = "'-1.23%"

# Row 22
$ws.Range("D22").Value = "'0.2616"
$ws.Range("E22").Value = "'5.05%"

# Row 24
$ws.Range("D24").Value = "'0.04371"
$ws.Range("E24").Value = "'1.26%"

# Row 25
$ws.Range("D25").Value = "'0.001239"
$ws.Range("E25").Value = "'1.40%"

# Row 26
$ws.Range("D26").Value = "'0.004303"
$ws.Range("E26").Value = "'-5.38%"

# Row 27
$ws.Range("D27").Value = "'0.0001220"
$ws.Range("E27").Value = "'-9.83%"

# Row 39
$ws.Range("D39").Value = "'0.02124"
$ws.Range("E39").Value = "'-6.60%"

# Row 40
$ws.Range("D40").Value = "'0.05173"
$ws.Range("E40").Value = "'2.43%"

# Row 41
$ws.Range("D41").Value = "'0.007563"
$ws.Range("E41").Value = "'1.88%"

# Row 42
$ws.Range("D42").Value = "'0.009137"
$ws.Range("E42").Value = "'-7.91%"

# Row 43
$ws.Range("D43").Value = "'0.1361"
$ws.Range("E43").Value = "'0.62%"

# Row 44
$ws.Range("D44").Value = "'0.002050"
$ws.Range("E44").Value = "'0.27%"

# Row 45
$ws.Range("D45").Value = "'0.008616"
$ws.Range("E45").Value = "'-1.82%"

# Row 46
$ws.Range("D46").Value = "'0.00006690"
$ws.Range("E46").Value = "'-1.07%"

# Row 47
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.19%"

# Row 48
$ws.Range("D48").Value = "'0.001200"
$ws.Range("E48").Value = "'-7.80%"

# Row 49
$ws.Range("D49").Value = "'0.002839"
$ws.Range("E49").Value = "'-5.53%"

# Row 50
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.19%"

# Row 51
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.19%"
